$wb = $excel.ActiveWorkbook

# --- Part 1: add the "2022-Q4" sheet, positioned right after "总计" ---
# Copy "2022-Q3" (same column layout/styling as every quarterly detail sheet)
# to just before itself, then rename the copy and rewrite its data for Q4.
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3, $null)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# The source sheet has 11 data rows (rows 2-12); Q4 only needs 6 (rows 2-7).
$q4.Range("A8:H12").EntireRow.Delete()

# Columns B-G hold fund code / name / size / position values. Several look
# numeric ("006511", "0.0560", ...) but must stay text so leading/trailing
# zeros survive - force the Text number format before assigning them.
$q4.Range("B2:G7").NumberFormat = "@"

$q4Data = @(
    @(0, "006511", "博道卓远混合A",                 "3.12", "81.13", "2.23", "0.0696", 8),
    @(1, "014663", "富国创新发展两年定期开放混合A", "2.26", "82.02", "2.48", "0.0560", 6),
    @(2, "007826", "博道志远混合C",                 "1.57", "82.17", "2.54", "0.0399", 5),
    @(3, "007825", "博道志远混合A",                 "0.99", "82.17", "2.54", "0.0251", 5),
    @(4, "006512", "博道卓远混合C",                 "1.02", "81.13", "2.23", "0.0227", 8),
    @(5, "014664", "富国创新发展两年定期开放混合C", "0.33", "82.02", "2.48", "0.0082", 6)
)

$r = 2
foreach ($row in $q4Data) {
    $q4.Cells.Item($r, 1).Value = $row[0]
    $q4.Cells.Item($r, 2).Value = $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = $row[3]
    $q4.Cells.Item($r, 5).Value = $row[4]
    $q4.Cells.Item($r, 6).Value = $row[5]
    $q4.Cells.Item($r, 7).Value = $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# --- Part 2: update "总计" summary sheet - insert the Q4 row, renumber the rest ---
$total = $wb.Worksheets.Item("总计")
$total.Range("A2").EntireRow.Insert()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 6
$total.Cells.Item(2, 4).Value = 0.22

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3

# Carry the index-column style (bold + border) onto the new A2 cell, and
# strip the stray formatting the row-insert left behind on B2:D2.
$total.Cells.Item(3, 1).Copy()
$total.Cells.Item(2, 1).PasteSpecial(-4122)
$total.Range("B2:D2").ClearFormats()

# Restore the originally-active sheet/tab (copying "2022-Q3" made the new
# "2022-Q4" copy active instead).
$wb.Worksheets.Item("2022-Q1").Activate()
